$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 39 (everything below shifts down by one,
# so the old row 39 becomes row 40, ..., and the old last row 180 becomes the new row 181).
$ws.Rows.Item(39).Insert()

# Populate the new row 39 with the new record (same fixed columns as the
# neighbouring rows, new Fecha/Volumen/Precio values).
$ws.Cells.Item(39, 1).Value = 5
$ws.Cells.Item(39, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(39, 3).Value = "Maule"
$ws.Cells.Item(39, 4).Value = 44715
$ws.Cells.Item(39, 5).Value = 7
$ws.Cells.Item(39, 6).Value = 100112017
$ws.Cells.Item(39, 7).Value = "Apio"
$ws.Cells.Item(39, 8).Value = "Americana (o)"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 700
$ws.Cells.Item(39, 11).Value = 6000
$ws.Cells.Item(39, 12).Value = 6000
$ws.Cells.Item(39, 13).Value = 6000
$ws.Cells.Item(39, 14).Value = "$/docena de matas"
$ws.Cells.Item(39, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(39, 16).Value = 1000
$ws.Cells.Item(39, 17).Value = 6
$ws.Cells.Item(39, 18).Value = "Hortaliza"
